# Insert a new weekly record as row 75, shifting the existing rows
# 75-98 down to 76-99 (Arveja Verde / Macroferia Regional de Talca sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value = 5
$ws.Range("B75").Value = "Macroferia Regional de Talca"
$ws.Range("C75").Value = "Maule"
$ws.Range("D75").Value = 44642
$ws.Range("E75").Value = 7
$ws.Range("F75").Value = 100112022
$ws.Range("G75").Value = "Arveja Verde"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 200
$ws.Range("K75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("M75").Value = 25000
$ws.Range("N75").Value = "`$/saco 25 kilos"
$ws.Range("O75").Value = "Carahue"
$ws.Range("P75").Value = 1000
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = "Hortaliza"
